$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row, taken from the diff's "after" (+) side.
# Columns: D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
#          M (Precio promedio ponderado), P (Precio $/Kg)

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13)
$D = @{ 2=44497; 3=44484; 4=44447; 5=44496; 6=44476; 7=44453; 8=44203; 9=44474; 10=44452; 11=44473; 12=44483; 13=44487 }
$J = @{ 2=50;    3=40;    4=75;    5=40;    6=30;    7=20;    8=30;    9=20;    10=120;   11=140;   12=50;    13=50 }
$K = @{ 2=2200;  3=2200;  4=2200;  5=2200;  6=2200;  7=2300;  8=2000;  9=1600;  10=2300;  11=1600;  12=2200;  13=2200 }
$L = @{ 2=2200;  3=2200;  4=2200;  5=2200;  6=2200;  7=2300;  8=2000;  9=1600;  10=2300;  11=1600;  12=2200;  13=2200 }
$M = @{ 2=2200;  3=2200;  4=2200;  5=2200;  6=2200;  7=2300;  8=2000;  9=1600;  10=2300;  11=1600;  12=2200;  13=2200 }
$P = @{ 2=2200;  3=2200;  4=2200;  5=2200;  6=2200;  7=2300;  8=2000;  9=1600;  10=2300;  11=1600;  12=2200;  13=2200 }

foreach ($r in $rows) {
    $ws.Range("D$r").Value = $D[$r]
    $ws.Range("J$r").Value = $J[$r]
    $ws.Range("K$r").Value = $K[$r]
    $ws.Range("L$r").Value = $L[$r]
    $ws.Range("M$r").Value = $M[$r]
    $ws.Range("P$r").Value = $P[$r]
}
